# Moved statistical datasets and results
#
# The sheet gains a new leading label column and a new header row:
#   1) Insert a new column before column A (existing A:D data shifts to B:E)
#   2) Insert a new row before row 1 (existing rows 1:21 shift to 2:22)
#   3) Fill the new header row (row 1, columns B:E) with column titles
#   4) Fill the new label column (A2:A22) with descriptive row names
#   5) Size the new label column to fit its (long) text

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("A:A").Insert()
$ws.Rows("1:1").Insert()

# New header row
$ws.Cells.Item(1, 2).Value = "Valid"
$ws.Cells.Item(1, 3).Value = "T"
$ws.Cells.Item(1, 4).Value = "Z"
$ws.Cells.Item(1, 5).Value = "p-value"

# New row labels (column A), one per former data row
$labels = @(
    "CyclomaticComplexity(CC) & CyclomaticComplexity(CC)",
    "MaintainabilityIndex & MaintainabilityIndex",
    "NbUniqueOperands & NbUniqueOperands",
    "NbUniqueOperands & EffortToImplement",
    "NbOperands & NbOperands",
    "NbOperands & EffortToImplement",
    "NbUniqueOperators & NbUniqueOperators",
    "NbUniqueOperators & EffortToImplement",
    "NbOperators & NbOperators",
    "ProgramLength & ProgramLength",
    "ProgramLength & EffortToImplement",
    "VocabularySize & VocabularySize",
    "ProgramVolume & ProgramVolume",
    "DifficultyLevel & DifficultyLevel",
    "ProgramLevel & ProgramLevel",
    "EffortToImplement & NbUniqueOperands",
    "EffortToImplement & NbOperands",
    "EffortToImplement & NbUniqueOperators",
    "EffortToImplement & ProgramLength",
    "EffortToImplement & EffortToImplement",
    "TimeToImplement & TimeToImplement"
)

for ($i = 0; $i -lt $labels.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $labels[$i]
}

# Widen the new label column to fit the long descriptive text.
# (Columns B:E keep the widths they inherited from the insert, untouched.)
$ws.Columns.Item(1).ColumnWidth = 53.666666666666664

Write-Output "ok"
